# Added check for empty rows.
#
# Summary of changes applied to the workbook:
#  1. Profile sheet (sheet1): append an empty "quote-prefixed" cell at A7
#     (dimension grows to A1:J7) and move the selection to A5.
#  2. Add a new "Expenses" sheet (sheetId 3) after Product, and make it the
#     active sheet (activeTab moves to index 2; Product loses tabSelected).
#  3. Populate the Expenses sheet with an entry list that intentionally
#     skips row 8, to exercise "empty row" handling: rows 1-2 are reserved
#     rows (italic), row 3 is a bold header (Entry No / Amount), rows 4-7
#     and 9 hold entries A001/A2/A3/A4/A6 with amounts, row 8 is left empty.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Profile sheet: add the trailing empty row used as an "empty row" probe
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Typing a lone apostrophe stores an empty, text-quote-prefixed cell (the
# xlsx <c s=".." /> with quotePrefix="1" on its style, and no value/type).
$ws1.Range("A7").Value = "'"
$ws1.Range("A7").Value = ""

# Move the active selection to A5, like the authored workbook.
$ws1.Range("A5").Select()

# ---------------------------------------------------------------------
# 2) Add the new "Expenses" worksheet after the existing sheets
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "Expenses"

# ---------------------------------------------------------------------
# 3) Fill in the Expenses sheet content
# ---------------------------------------------------------------------
# Rows 1-2: reserved-row marker, italic (same style as Product!A1/A2).
$ws3.Range("A1").Value = "Reserved row"
$ws3.Range("A1").Font.Italic = $true
$ws3.Range("A2").Value = "Reserved row"
$ws3.Range("A2").Font.Italic = $true

# Row 3: bold header row.
$ws3.Range("B3").Value = "Amount"
$ws3.Range("A3").Value = "Entry No"
$ws3.Range("A3:B3").Font.Bold = $true

# Rows 4-7: entries A001, A2, A3, A4 with their amounts. The "Entry No"
# column uses a custom date-ish number format (d-mmm) even though it holds
# text, matching the authored workbook's styling.
$ws3.Range("A4").NumberFormat = "d-mmm"
$ws3.Range("A4").Value = "A001"
$ws3.Range("B4").Value = 3.09

$ws3.Range("A5").NumberFormat = "d-mmm"
$ws3.Range("A5").Value = "A2"
$ws3.Range("B5").Value = 2.55

$ws3.Range("A6").NumberFormat = "d-mmm"
$ws3.Range("A6").Value = "A3"
$ws3.Range("B6").Value = 10.77

$ws3.Range("A7").NumberFormat = "d-mmm"
$ws3.Range("A7").Value = "A4"
$ws3.Range("B7").Value = 5.6

# Row 8 intentionally left blank (the "empty row" being checked for).

# Row 9: one more entry, A6, after the gap.
$ws3.Range("A9").NumberFormat = "d-mmm"
$ws3.Range("A9").Value = "A6"
$ws3.Range("B9").Value = 2.56

# Selection on the new sheet sits just past the data, at A10.
$ws3.Range("A10").Select()
